$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (rows 2-21) down by 4 rows (to rows 6-25),
# copying values cell-by-cell (bottom-up so sources aren't clobbered before
# being read). This avoids Insert()'s side effect of pulling in the header
# row's bold style into the newly created rows.
for ($r = 21; $r -ge 2; $r--) {
    $destRow = $r + 4
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value2
}

# Fill the 4 freshly freed-up rows at the top (rows 2-5) with the new data.
$topData = @(
    @(-0.0056505035609006, -0.007177666760981, 0),
    @(-0.0161879286170005, 0.0122173046693205, -0.0047342055477201),
    @(0.0029016099870204, -0.0010690141934901, -0.009468411095440299),
    @(0.00534507073462, 0.0088575463742017, 0.0045814891345798)
)

for ($i = 0; $i -lt $topData.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $topData[$i][0]
    $ws.Cells.Item($row, 2).Value = $topData[$i][1]
    $ws.Cells.Item($row, 3).Value = $topData[$i][2]
}

# Append 6 brand new rows at the bottom (rows 26-31).
$bottomData = @(
    @(-0.6475171446800232, -0.2813034355640411, 0.0232128798961639),
    @(-0.1372919678688049, -2.705674886703491, -0.5198463201522827),
    @(-0.6478226184844971, -0.2128865420818328, -0.0656680166721344),
    @(-0.1000291854143142, 0.1372919678688049, -0.1838704347610473),
    @(0.2654209434986114, 0.0520762614905834, 0.0438295826315879),
    @(0.0809396430850029, 0.3397938013076782, 0.0826195254921913)
)

for ($i = 0; $i -lt $bottomData.Length; $i++) {
    $row = $i + 26
    $ws.Cells.Item($row, 1).Value = $bottomData[$i][0]
    $ws.Cells.Item($row, 2).Value = $bottomData[$i][1]
    $ws.Cells.Item($row, 3).Value = $bottomData[$i][2]
}
